# Update the handoff/handback timestamps for the 30609edd... entry
# (the one that generated the report) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-15 03:36:06"
$wsZhCn.Range("G4").Value = "2016-02-15 03:36:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-15 03:36:20"
$wsDeDe.Range("G4").Value = "2016-02-15 03:37:23"
